# feat: Published NHS Equity and Skin Tone standards
#
# 1) Fix the bad "null#nursing" category-code text (shared by the three
#    existing "braden scale assessment" / "goal evaluation" / "nursing
#    assessment" rows) so it reads "Observation Category Codes#nursing".
# 2) Append a new "ONC Skin Tone Observation" profile as row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the Category Code text everywhere it is used -------------
$ws.Range("C2").Value = "Observation Category Codes#nursing"
$ws.Range("C9").Value = "Observation Category Codes#nursing"
$ws.Range("C10").Value = "Observation Category Codes#nursing"

# --- 2. Add the new "ONC Skin Tone Observation" row -----------------------
$ws.Range("A11").Value = "onc-skintone-observation"
$ws.Range("B11").Value = "ONC Skin Tone Observation"
# Category Code / Category VS are blank for this profile; a leading quote
# forces a real (empty) text cell instead of clearing it outright.
$ws.Range("C11").Value = "'"
$ws.Range("D11").Value = "'"
$ws.Range("E11").Value = "LOINC#66472-2"
# Code VS is blank too.
$ws.Range("F11").Value = "'"
$ws.Range("G11").Value = "dateTime, Period, Timing, instant"
$ws.Range("H11").Value = "CodeableConceptĵ"
$ws.Range("I11").Value = "optional"
# Body Site / Method are blank.
$ws.Range("J11").Value = "'"
$ws.Range("K11").Value = "'"

# Match the existing table's row styling (border + top/wrap alignment) by
# copying the formatting from the row directly above onto the new row.
$ws.Range("A10:K10").Copy()
$ws.Range("A11:K11").PasteSpecial(-4122)
